# Apply updated crypto price/volume snapshot to Sheet1 (cryptos.xlsx).
# Numeric-looking Price values are written with a leading quote then
# reset to the "Normal" style so they round-trip as plain text cells
# (matching the original inline-string cells) instead of being coerced
# into Excel numbers (which would lose trailing zeros / exact digits).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "41.045.33"
$ws.Range("E2").Value = "  -2.15%  "

$ws.Range("D3").Value = "2.134.22"
$ws.Range("E3").Value = "  -3.72%  "

$ws.Range("E4").Value = "  +0.17%  "

Set-TextValue $ws.Range("D5") "234.69"
$ws.Range("E5").Value = "  -2.83%  "

Set-TextValue $ws.Range("D6") "0.596"
$ws.Range("E6").Value = "  -4.92%  "

Set-TextValue $ws.Range("D7") "68.80"
$ws.Range("E7").Value = "  -5.99%  "

$ws.Range("E8").Value = "  +0.15%  "

Set-TextValue $ws.Range("D9") "0.562"
$ws.Range("E9").Value = "  -7.11%  "

Set-TextValue $ws.Range("D10") "38.14"
$ws.Range("E10").Value = "  -10.19%  "

Set-TextValue $ws.Range("D11") "0.0885"
$ws.Range("E11").Value = "  -7.36%  "

Set-TextValue $ws.Range("D12") "53.24"
$ws.Range("E12").Value = "  -6.96%  "

Set-TextValue $ws.Range("D13") "0.0992"
$ws.Range("E13").Value = "  -4.31%  "

Set-TextValue $ws.Range("D14") "6.53"
$ws.Range("E14").Value = "  -6.62%  "

$ws.Range("D15").Value = "2.455.15"
$ws.Range("E15").Value = "  -3.66%  "

Set-TextValue $ws.Range("D16") "14.11"
$ws.Range("E16").Value = "  -0.88%  "

$ws.Range("D17").Value = "2.112.47"
$ws.Range("E17").Value = "  -4.46%  "

Set-TextValue $ws.Range("D18") "0.769"
$ws.Range("E18").Value = "  -7.75%  "

$ws.Range("D19").Value = "40.894.65"
$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  -7.21%  "

Set-TextValue $ws.Range("D21") "68.50"
$ws.Range("E21").Value = "  -5.80%  "

Set-TextValue $ws.Range("D22") "5.67"
$ws.Range("E22").Value = "  -8.15%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D23") "223.57"
$ws.Range("E23").Value = "  -2.73%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D24") "9.38"
$ws.Range("E24").Value = "  -13.33%  "

$ws.Range("E25").Value = "  -7.34%  "

$ws.Range("E26").Value = "  -0.24%  "

Set-TextValue $ws.Range("D27") "10.43"
$ws.Range("E27").Value = "  -9.87%  "

Set-TextValue $ws.Range("D28") "3.25"
$ws.Range("E28").Value = "  -13.28%  "

$ws.Range("E29").Value = "  -1.25%  "

Set-TextValue $ws.Range("D30") "2.13"
$ws.Range("E30").Value = "  -6.22%  "

Set-TextValue $ws.Range("D31") "169.89"
$ws.Range("E31").Value = "  +1.34%  "

Set-TextValue $ws.Range("D32") "19.44"
$ws.Range("E32").Value = "  -5.15%  "

Set-TextValue $ws.Range("D33") "30.70"
$ws.Range("E33").Value = "  +1.97%  "

Set-TextValue $ws.Range("D34") "0.0742"
$ws.Range("E34").Value = "  -6.87%  "

Set-TextValue $ws.Range("D35") "5.01"
$ws.Range("E35").Value = "  -11.41%  "

$ws.Range("E36").Value = "  -4.99%  "

Set-TextValue $ws.Range("D37") "0.102"
$ws.Range("E37").Value = "  -7.41%  "

Set-TextValue $ws.Range("D38") "4.09"
$ws.Range("E38").Value = "  -3.40%  "

Set-TextValue $ws.Range("D39") "0.0284"
$ws.Range("E39").Value = "  -5.39%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D40") "2.02"
$ws.Range("E40").Value = "  -4.73%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D41") "11.54"
$ws.Range("E41").Value = "  -16.26%  "

Set-TextValue $ws.Range("D42") "5.17"
$ws.Range("E42").Value = "  -8.31%  "

Set-TextValue $ws.Range("D43") "56.99"
$ws.Range("E43").Value = "  -12.26%  "

$ws.Range("E44").Value = "  -6.88%  "

$ws.Range("E45").Value = "  -7.22%  "

Set-TextValue $ws.Range("D46") "0.0951"
$ws.Range("E46").Value = "  -5.06%  "

Set-TextValue $ws.Range("D47") "96.00"
$ws.Range("E47").Value = "  -8.34%  "

Set-TextValue $ws.Range("D48") "1.06"
$ws.Range("E48").Value = "  -4.60%  "

$ws.Range("E49").Value = "  -5.86%  "

Set-TextValue $ws.Range("D50") "2.61"
$ws.Range("E50").Value = "  -3.23%  "

Set-TextValue $ws.Range("D51") "2.11"
$ws.Range("E51").Value = "  -11.19%  "
